$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly "Apio" (celery) price block for "Vega Monumental Concepción"
# occupies rows 197-244. A new week's pricing (rows for "Primera" and
# "Segunda" quality) was inserted at the top of that block (new rows 197
# and 198), pushing all the existing rows down by two. That also pushes
# the two rows that used to be at the bottom of the block (243, 244) past
# the old end of the sheet, so they land in two newly appended rows
# (245, 246).

$firstDataRow = 197
$lastDataRow = 244
$shift = 2
$newLastRow = $lastDataRow + $shift

# Shift existing rows 197..244 down to 199..246, working from the bottom
# up so a row is always read before anything gets written into it.
for ($destRow = $newLastRow; $destRow -ge ($firstDataRow + $shift); $destRow--) {
    $srcRow = $destRow - $shift

    for ($col = 1; $col -le 18; $col++) {
        $val = $ws.Cells.Item($srcRow, $col).Value2
        $ws.Cells.Item($destRow, $col).Value = $val
    }
}

# Populate the two new rows (197 = "Primera", 198 = "Segunda") with the
# new week's data. The descriptive columns are identical to every other
# row in this block.
$ws.Cells.Item(197, 1).Value = 11
$ws.Cells.Item(197, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(197, 3).Value = "Bíobío"
$ws.Cells.Item(197, 4).Value = 44663
$ws.Cells.Item(197, 5).Value = 8
$ws.Cells.Item(197, 6).Value = 100112017
$ws.Cells.Item(197, 7).Value = "Apio"
$ws.Cells.Item(197, 8).Value = "Americana (o)"
$ws.Cells.Item(197, 9).Value = "Primera"
$ws.Cells.Item(197, 10).Value = 100
$ws.Cells.Item(197, 11).Value = 7500
$ws.Cells.Item(197, 12).Value = 8000
$ws.Cells.Item(197, 13).Value = 7750
$ws.Cells.Item(197, 14).Value = "$/docena de matas"
$ws.Cells.Item(197, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(197, 16).Value = 1292
$ws.Cells.Item(197, 17).Value = 6
$ws.Cells.Item(197, 18).Value = "Hortaliza"

$ws.Cells.Item(198, 1).Value = 11
$ws.Cells.Item(198, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(198, 3).Value = "Bíobío"
$ws.Cells.Item(198, 4).Value = 44663
$ws.Cells.Item(198, 5).Value = 8
$ws.Cells.Item(198, 6).Value = 100112017
$ws.Cells.Item(198, 7).Value = "Apio"
$ws.Cells.Item(198, 8).Value = "Americana (o)"
$ws.Cells.Item(198, 9).Value = "Segunda"
$ws.Cells.Item(198, 10).Value = 50
$ws.Cells.Item(198, 11).Value = 6500
$ws.Cells.Item(198, 12).Value = 6500
$ws.Cells.Item(198, 13).Value = 6500
$ws.Cells.Item(198, 14).Value = "$/docena de matas"
$ws.Cells.Item(198, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(198, 16).Value = 1083
$ws.Cells.Item(198, 17).Value = 6
$ws.Cells.Item(198, 18).Value = "Hortaliza"

# Make sure every "Fecha" cell in the (now two-rows-longer) block keeps
# the date number format used throughout column D, including the two
# brand-new rows at the very bottom that the shift produced.
$ws.Range("D$firstDataRow`:D$newLastRow").NumberFormat = "YYYY-MM-DD HH:MM:SS"
